$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 0.8200883333333334
$ws.Range("H2").Value2 = 2.460265
$ws.Range("I2").Value2 = 0.2405117342909232
$ws.Range("J2").Value2 = 0.2405117342909232
$ws.Range("M2").Value2 = 0.01102566666666667
$ws.Range("N2").Value2 = 0.033077
$ws.Range("O2").Value2 = 0.07220475878629121
$ws.Range("P2").Value2 = 0.07220475878629119
$ws.Range("Q2").Value2 = 0.009042020600555558
$ws.Range("R2").Value2 = 0.08137818540500001
$ws.Range("S2").Value2 = 0.01736609175974867
$ws.Range("T2").Value2 = 0.01736609175974867

# Row 3
$ws.Range("G3").Value2 = 0.8200883333333334
$ws.Range("H3").Value2 = 2.460265
$ws.Range("I3").Value2 = 0.2405117342909232
$ws.Range("J3").Value2 = 0.2405117342909232
$ws.Range("O3").Value2 = 0.1778345339445536
$ws.Range("P3").Value2 = 0.1778345339445536
$ws.Range("Q3").Value2 = 0.02226977205444445
$ws.Range("R3").Value2 = 0.20042794849
$ws.Range("S3").Value2 = 0.04277129217582263
$ws.Range("T3").Value2 = 0.04277129217582262

# Row 4
$ws.Range("G4").Value2 = 0.8200883333333334
$ws.Range("H4").Value2 = 2.460265
$ws.Range("I4").Value2 = 0.2405117342909232
$ws.Range("J4").Value2 = 0.2405117342909232
$ws.Range("M4").Value2 = 0.07542833333333333
$ws.Range("N4").Value2 = 0.226285
$ws.Range("O4").Value2 = 0.4939641999563414
$ws.Range("P4").Value2 = 0.4939641999563414
$ws.Range("Q4").Value2 = 0.06185789616944445
$ws.Range("R4").Value2 = 0.5567210655250001
$ws.Range("S4").Value2 = 0.118804186409128
$ws.Range("T4").Value2 = 0.118804186409128

# Row 5
$ws.Range("G5").Value2 = 0.8200883333333334
$ws.Range("H5").Value2 = 2.460265
$ws.Range("I5").Value2 = 0.2405117342909232
$ws.Range("J5").Value2 = 0.2405117342909232
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.03909066666666667
$ws.Range("N5").Value2 = 0.117272
$ws.Range("O5").Value2 = 0.2559965073128138
$ws.Range("P5").Value2 = 0.2559965073128138
$ws.Range("Q5").Value2 = 0.03205779967555556
$ws.Range("R5").Value2 = 0.28852019708
$ws.Range("S5").Value2 = 0.06157016394622384
$ws.Range("T5").Value2 = 0.06157016394622384

# Row 6
$ws.Range("A6").Value2 = "Resolving-Mac"
$ws.Range("D6").Value2 = "ECs"
$ws.Range("G6").Value2 = 2.589676
$ws.Range("H6").Value2 = 7.769028
$ws.Range("I6").Value2 = 0.7594882657090768
$ws.Range("J6").Value2 = 0.7594882657090768
$ws.Range("M6").Value2 = 0.01102566666666667
$ws.Range("N6").Value2 = 0.033077
$ws.Range("O6").Value2 = 0.07220475878629121
$ws.Range("P6").Value2 = 0.07220475878629119
$ws.Range("Q6").Value2 = 0.02855290435066667
$ws.Range("R6").Value2 = 0.256976139156
$ws.Range("S6").Value2 = 0.05483866702654253
$ws.Range("T6").Value2 = 0.05483866702654253

# Row 7
$ws.Range("D7").Value2 = "FAPs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 2.589676
$ws.Range("H7").Value2 = 7.769028
$ws.Range("I7").Value2 = 0.7594882657090768
$ws.Range("J7").Value2 = 0.7594882657090768
$ws.Range("K7").Value2 = 2
$ws.Range("L7").Value2 = 0.6666666666666666
$ws.Range("M7").Value2 = 0.02715533333333333
$ws.Range("N7").Value2 = 0.081466
$ws.Range("O7").Value2 = 0.1778345339445536
$ws.Range("P7").Value2 = 0.1778345339445536
$ws.Range("Q7").Value2 = 0.07032351500533333
$ws.Range("R7").Value2 = 0.632911635048
$ws.Range("S7").Value2 = 0.1350632417687309
$ws.Range("T7").Value2 = 0.1350632417687309

# Row 8
$ws.Range("D8").Value2 = "Inflammatory-Mac"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 2.589676
$ws.Range("H8").Value2 = 7.769028
$ws.Range("I8").Value2 = 0.7594882657090768
$ws.Range("J8").Value2 = 0.7594882657090768
$ws.Range("M8").Value2 = 0.07542833333333333
$ws.Range("N8").Value2 = 0.226285
$ws.Range("O8").Value2 = 0.4939641999563414
$ws.Range("P8").Value2 = 0.4939641999563414
$ws.Range("Q8").Value2 = 0.1953349445533333
$ws.Range("R8").Value2 = 1.75801450098
$ws.Range("S8").Value2 = 0.3751600135472133
$ws.Range("T8").Value2 = 0.3751600135472133

# Row 9
$ws.Range("D9").Value2 = "MuSCs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 2.589676
$ws.Range("H9").Value2 = 7.769028
$ws.Range("I9").Value2 = 0.7594882657090768
$ws.Range("J9").Value2 = 0.7594882657090768
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 0.03909066666666667
$ws.Range("N9").Value2 = 0.117272
$ws.Range("O9").Value2 = 0.2559965073128138
$ws.Range("P9").Value2 = 0.2559965073128138
$ws.Range("Q9").Value2 = 0.1012321612906667
$ws.Range("R9").Value2 = 0.9110894516160001
$ws.Range("S9").Value2 = 0.19442634336659
$ws.Range("T9").Value2 = 0.19442634336659

# Remove the two trailing rows (old rows 10 and 11) that no longer exist
$ws.Rows("10:11").Delete()

Write-Host "edit applied"
